# Refresh the "cryptos" price/volume table (GitHub Actions daily update).
# All Price (D) and Volume(1h) (E) values are stored as text, and two rows
# (Celestia / ARBITRUM) also swap position in the ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving purely numeric-looking text need to be forced to Text format
# before assignment (otherwise Excel auto-converts "3.60" -> 3.6), then have their
# format cleared again afterwards so the saved style matches the original (no explicit style).
$textCells = @("D5", "D6", "D10", "D11", "D12", "D15", "D16", "D18", "D21", "D22", "D23", "D24", "D26", "D30", "D31", "D32", "D34", "D36", "D38", "D39", "D40", "D41", "D42", "D44", "D46", "D47", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "38.871.62"
$ws.Range("E2").Value = "  -4.81%  "
$ws.Range("D3").Value = "2.224.36"
$ws.Range("E3").Value = "  -6.80%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "296.38"
$ws.Range("E5").Value = "  -5.70%  "
$ws.Range("D6").Value = "79.92"
$ws.Range("E6").Value = "  -9.65%  "
$ws.Range("E7").Value = "  -4.81%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -7.47%  "
$ws.Range("D10").Value = "0.0770"
$ws.Range("E10").Value = "  -6.68%  "
$ws.Range("D11").Value = "27.84"
$ws.Range("E11").Value = "  -10.98%  "
$ws.Range("D12").Value = "45.92"
$ws.Range("E12").Value = "  -13.56%  "
$ws.Range("D14").Value = "2.566.39"
$ws.Range("E14").Value = "  -6.79%  "
$ws.Range("D15").Value = "6.09"
$ws.Range("E15").Value = "  -7.53%  "
$ws.Range("D16").Value = "13.98"
$ws.Range("E16").Value = "  -7.73%  "
$ws.Range("D17").Value = "2.235.47"
$ws.Range("E17").Value = "  -5.14%  "
$ws.Range("D18").Value = "0.712"
$ws.Range("E18").Value = "  -6.65%  "
$ws.Range("D19").Value = "38.817.32"
$ws.Range("E19").Value = "  -4.72%  "
$ws.Range("D20").Value = "0.0₃0857"
$ws.Range("E20").Value = "  -6.25%  "
$ws.Range("D21").Value = "5.71"
$ws.Range("E21").Value = "  -7.87%  "
$ws.Range("D22").Value = "64.79"
$ws.Range("E22").Value = "  -6.61%  "
$ws.Range("D23").Value = "9.83"
$ws.Range("E23").Value = "  -9.21%  "
$ws.Range("D24").Value = "224.54"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "2.38"
$ws.Range("E26").Value = "  -10.37%  "
$ws.Range("E27").Value = "  -6.29%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("E29").Value = "  -7.00%  "
$ws.Range("D30").Value = "8.85"
$ws.Range("E30").Value = "  -5.84%  "
$ws.Range("D31").Value = "148.21"
$ws.Range("E31").Value = "  -5.31%  "
$ws.Range("D32").Value = "31.04"
$ws.Range("E32").Value = "  -8.81%  "
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("D34").Value = "4.76"
$ws.Range("E34").Value = "  -8.92%  "
$ws.Range("E35").Value = "  -3.93%  "
$ws.Range("D36").Value = "0.0683"
$ws.Range("E36").Value = "  -6.89%  "
$ws.Range("E37").Value = "  -4.56%  "
$ws.Range("D38").Value = "2.65"
$ws.Range("E38").Value = "  -5.99%  "
$ws.Range("D39").Value = "0.0954"
$ws.Range("E39").Value = "  -5.00%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.59"
$ws.Range("E40").Value = "  -8.52%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").Value = "14.37"
$ws.Range("E41").Value = "  -10.80%  "
$ws.Range("D42").Value = "3.60"
$ws.Range("E42").Value = "  -5.94%  "
$ws.Range("D43").Value = "1.895.90"
$ws.Range("E43").Value = "  -3.35%  "
$ws.Range("D44").Value = "2.05"
$ws.Range("E44").Value = "  -9.82%  "
$ws.Range("E45").Value = "  -7.08%  "
$ws.Range("D46").Value = "16.07"
$ws.Range("E46").Value = "  -8.95%  "
$ws.Range("D47").Value = "8.97"
$ws.Range("E47").Value = "  -4.11%  "
$ws.Range("E48").Value = "  -10.90%  "
$ws.Range("D49").Value = "2.438.44"
$ws.Range("E49").Value = "  -6.78%  "
$ws.Range("D50").Value = "68.29"
$ws.Range("E50").Value = "  -6.45%  "
$ws.Range("D51").Value = "87.23"
$ws.Range("E51").Value = "  -7.23%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
